# Update results for Steel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iron & steel / Hydrogen (B3) updated to new computed result
$ws.Range("B3").Value = 202.5141002774836

# Non-metallic minerals / Biomass (D6): refreshed floating point result
$ws.Range("D6").Value = 954.0503059973877
